# "add test skip lines"
# Insert a blank data row above row 6 (shifting the existing rows 6-8 down to
# 7-9), turning what used to be a contiguous data block into one that has a
# skipped/blank line in the middle - matching the commit's "add test skip
# lines" intent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above row 6, pushing rows 6-8 down to 7-9 -------------
$ws.Rows.Item(6).Insert()

# The freshly inserted row inherits generic formatting; copy the (identical)
# number/style formats from the row above (row 5) so the new blank row 6
# keeps using the same style ids as the rest of the data rows.
$ws.Range("A5:I5").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Keep the new blank row 6 at the same height as its neighbours.
$ws.Rows.Item(6).RowHeight = 13.8

# Row 8 (originally row 7, "Pau d'alho" / 42499 line) gets an explicit,
# slightly smaller custom row height after the shift.
$ws.Rows.Item(8).RowHeight = 12.8

# --- Update the worksheet view / selection ----------------------------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("B8").Select()

# --- Extend the (hidden) _FilterDatabase defined name by one row -----------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Plan1!`$A`$1:`$BJ`$207"
    }
}
